$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

function Set-TextValue($addr, $text) {
    $target = $ws.Range($addr)
    $escaped = $text -replace '"', '""'
    $target.Formula = '="' + $escaped + '"'
    $target.Copy()
    $target.PasteSpecial(-4163) # xlPasteValues
}

Set-TextValue "Q2"  "51524446"
Set-TextValue "Q3"  "51524447"
Set-TextValue "R3"  "12-18-2021"
Set-TextValue "AD3" "51524448"
Set-TextValue "Q4"  "51524449"

$excel.CutCopyMode = 0
